$d = $word.ActiveDocument

# Fix the temperature sensor list: "...,T9,T1,T1,pumpState..." -> "...,T9,T10,pumpState..."
$d.Content.Find.Execute(
    "T1,T2,T3,T4,T5,T6,T7,T8,T9,T1,T1,pumpState,pumpRpm,",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "T1,T2,T3,T4,T5,T6,T7,T8,T9,T10,pumpState,pumpRpm,",
    2
)
